# Refactoring and reading current processes from local machine with tested convention
#
# Add a new data row (row 2) underneath the existing header row, populating
# every column (A-H) with the value "mike" (this introduces a new shared
# string and leaves the header row/styles untouched), then leave the
# selection on C3 as the last active cell, matching the authored workbook
# state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:H2").Value = "mike"

$ws.Range("C3").Select()
